$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "B2m"
$ws.Range("C2").Value = "Gm11127"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1111.404703666667
$ws.Range("H2").Value = 3334.214111
$ws.Range("I2").Value = 0.2049713233760527
$ws.Range("J2").Value = 0.2049713233760527
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.185824
$ws.Range("N2").Value = 0.557472
$ws.Range("O2").Value = 0.7926810887363958
$ws.Range("P2").Value = 0.7926810887363958
$ws.Range("Q2").Value = 206.5256676541547
$ws.Range("R2").Value = 1858.731008887392
$ws.Range("S2").Value = 0.1624768917734693
$ws.Range("T2").Value = 0.1624768917734693

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "B2m"
$ws.Range("C3").Value = "Gm11127"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1111.404703666667
$ws.Range("H3").Value = 3334.214111
$ws.Range("I3").Value = 0.2049713233760527
$ws.Range("J3").Value = 0.2049713233760527
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.03181266666666666
$ws.Range("N3").Value = 0.095438
$ws.Range("O3").Value = 0.1357052869862955
$ws.Range("P3").Value = 0.1357052869862955
$ws.Range("Q3").Value = 35.35674736951311
$ws.Range("R3").Value = 318.210726325618
$ws.Range("S3").Value = 0.02781569226270802
$ws.Range("T3").Value = 0.02781569226270801

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "B2m"
$ws.Range("C4").Value = "Gm11127"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1111.404703666667
$ws.Range("H4").Value = 3334.214111
$ws.Range("I4").Value = 0.2049713233760527
$ws.Range("J4").Value = 0.2049713233760527
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.016788
$ws.Range("N4").Value = 0.050364
$ws.Range("O4").Value = 0.07161362427730872
$ws.Range("P4").Value = 0.0716136242773087
$ws.Range("Q4").Value = 18.658262165156
$ws.Range("R4").Value = 167.924359486404
$ws.Range("S4").Value = 0.01467873933987538
$ws.Range("T4").Value = 0.01467873933987538

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "B2m"
$ws.Range("C5").Value = "Gm11127"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1149.586873333333
$ws.Range("H5").Value = 3448.76062
$ws.Range("I5").Value = 0.2120130875688133
$ws.Range("J5").Value = 0.2120130875688133
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.185824
$ws.Range("N5").Value = 0.557472
$ws.Range("O5").Value = 0.7926810887363958
$ws.Range("P5").Value = 0.7926810887363958
$ws.Range("Q5").Value = 213.6208311502933
$ws.Range("R5").Value = 1922.58748035264
$ws.Range("S5").Value = 0.1680587650804117
$ws.Range("T5").Value = 0.1680587650804117

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "B2m"
$ws.Range("C6").Value = "Gm11127"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1149.586873333333
$ws.Range("H6").Value = 3448.76062
$ws.Range("I6").Value = 0.2120130875688133
$ws.Range("J6").Value = 0.2120130875688133
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.03181266666666666
$ws.Range("N6").Value = 0.095438
$ws.Range("O6").Value = 0.1357052869862955
$ws.Range("P6").Value = 0.1357052869862955
$ws.Range("Q6").Value = 36.57142400572889
$ws.Range("R6").Value = 329.14281605156
$ws.Range("S6").Value = 0.02877129689337641
$ws.Range("T6").Value = 0.02877129689337641

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "B2m"
$ws.Range("C7").Value = "Gm11127"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1149.586873333333
$ws.Range("H7").Value = 3448.76062
$ws.Range("I7").Value = 0.2120130875688133
$ws.Range("J7").Value = 0.2120130875688133
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.016788
$ws.Range("N7").Value = 0.050364
$ws.Range("O7").Value = 0.07161362427730872
$ws.Range("P7").Value = 0.0716136242773087
$ws.Range("Q7").Value = 19.29926442952
$ws.Range("R7").Value = 173.69337986568
$ws.Range("S7").Value = 0.01518302559502515
$ws.Range("T7").Value = 0.01518302559502514

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "B2m"
$ws.Range("C8").Value = "Gm11127"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2939.417277333333
$ws.Range("H8").Value = 8818.251832
$ws.Range("I8").Value = 0.5421033825947782
$ws.Range("J8").Value = 0.5421033825947781
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.185824
$ws.Range("N8").Value = 0.557472
$ws.Range("O8").Value = 0.7926810887363958
$ws.Range("P8").Value = 0.7926810887363958
$ws.Range("Q8").Value = 546.2142761431893
$ws.Range("R8").Value = 4915.928485288704
$ws.Range("S8").Value = 0.4297150995229116
$ws.Range("T8").Value = 0.4297150995229116

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "B2m"
$ws.Range("C9").Value = "Gm11127"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2939.417277333333
$ws.Range("H9").Value = 8818.251832
$ws.Range("I9").Value = 0.5421033825947782
$ws.Range("J9").Value = 0.5421033825947781
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.03181266666666666
$ws.Range("N9").Value = 0.095438
$ws.Range("O9").Value = 0.1357052869862955
$ws.Range("P9").Value = 0.1357052869862955
$ws.Range("Q9").Value = 93.51070203804622
$ws.Range("R9").Value = 841.596318342416
$ws.Range("S9").Value = 0.07356629511126593
$ws.Range("T9").Value = 0.07356629511126592

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "B2m"
$ws.Range("C10").Value = "Gm11127"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2939.417277333333
$ws.Range("H10").Value = 8818.251832
$ws.Range("I10").Value = 0.5421033825947782
$ws.Range("J10").Value = 0.5421033825947781
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.016788
$ws.Range("N10").Value = 0.050364
$ws.Range("O10").Value = 0.07161362427730872
$ws.Range("P10").Value = 0.0716136242773087
$ws.Range("Q10").Value = 49.34693725187201
$ws.Range("R10").Value = 444.122435266848
$ws.Range("S10").Value = 0.03882198796060059
$ws.Range("T10").Value = 0.03882198796060057

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "B2m"
$ws.Range("C11").Value = "Gm11127"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 221.836001
$ws.Range("H11").Value = 665.508003
$ws.Range("I11").Value = 0.04091220646035591
$ws.Range("J11").Value = 0.04091220646035591
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.185824
$ws.Range("N11").Value = 0.557472
$ws.Range("O11").Value = 0.7926810887363958
$ws.Range("P11").Value = 0.7926810887363958
$ws.Range("Q11").Value = 41.222453049824
$ws.Range("R11").Value = 371.002077448416
$ws.Range("S11").Value = 0.03243033235960313
$ws.Range("T11").Value = 0.03243033235960313

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "B2m"
$ws.Range("C12").Value = "Gm11127"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 221.836001
$ws.Range("H12").Value = 665.508003
$ws.Range("I12").Value = 0.04091220646035591
$ws.Range("J12").Value = 0.04091220646035591
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.03181266666666666
$ws.Range("N12").Value = 0.095438
$ws.Range("O12").Value = 0.1357052869862955
$ws.Range("P12").Value = 0.1357052869862955
$ws.Range("Q12").Value = 7.057194754479333
$ws.Range("R12").Value = 63.514752790314
$ws.Range("S12").Value = 0.005552002718945173
$ws.Range("T12").Value = 0.005552002718945173

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "B2m"
$ws.Range("C13").Value = "Gm11127"
$ws.Range("D13").Value = "M2"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 221.836001
$ws.Range("H13").Value = 665.508003
$ws.Range("I13").Value = 0.04091220646035591
$ws.Range("J13").Value = 0.04091220646035591
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.016788
$ws.Range("N13").Value = 0.050364
$ws.Range("O13").Value = 0.07161362427730872
$ws.Range("P13").Value = 0.0716136242773087
$ws.Range("Q13").Value = 3.724182784788
$ws.Range("R13").Value = 33.517645063092
$ws.Range("S13").Value = 0.002929871381807611
$ws.Range("T13").Value = 0.00292987138180761
